$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and row 23 had their stimulus images (and associated stats in
# columns L:V) swapped after the image folder was cleaned up.
$cols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V")

foreach ($col in $cols) {
    $cellA = $ws.Range("$col" + "2")
    $cellB = $ws.Range("$col" + "23")

    $valA = $cellA.Value2
    $valB = $cellB.Value2

    $cellA.Value2 = $valB
    $cellB.Value2 = $valA
}
